$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.476.27"
$ws.Range("D3").Value = "2.437.05"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "2.432.72"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "2.875.91"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "62.296.09"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "2.435.86"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.41%  "
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.48"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "632.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.08"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.41%  "
$ws.Range("D28").Value = "0.0₃0960"
$ws.Range("E28").Value = "  -6.13%  "
$ws.Range("D29").Value = "2.558.75"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  -3.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -5.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.48"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.23"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.12"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("E45").Value = "  -5.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.21"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.70"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.73"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.19%  "
$ws.Range("E51").Value = "  -1.99%  "
